$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.900.22"
$ws.Range("E2").Value = "  +2.05%  "

# Row 3
$ws.Range("D3").Value = "3.453.91"
$ws.Range("E3").Value = "  +1.61%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'576.14"
$ws.Range("E5").Value = "  +1.45%  "

# Row 6
$ws.Range("D6").Value = "'160.43"
$ws.Range("E6").Value = "  +3.07%  "

# Row 7
$ws.Range("D7").Value = "'0.609"
$ws.Range("E7").Value = "  +6.76%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").Value = "3.456.05"
$ws.Range("E9").Value = "  +1.59%  "

# Row 11
$ws.Range("E11").Value = "  +1.88%  "

# Row 12
$ws.Range("E12").Value = "  +3.36%  "

# Row 13
$ws.Range("D13").Value = "4.047.59"
$ws.Range("E13").Value = "  +1.42%  "

# Row 14
$ws.Range("E14").Value = "  +0.79%  "

# Row 15
$ws.Range("E15").Value = "  +1.44%  "

# Row 16
$ws.Range("D16").Value = "'28.26"
$ws.Range("E16").Value = "  +2.97%  "

# Row 17
$ws.Range("D17").Value = "64.922.47"
$ws.Range("E17").Value = "  +2.06%  "

# Row 18
$ws.Range("D18").Value = "3.493.87"
$ws.Range("E18").Value = "  +3.13%  "

# Row 19
$ws.Range("E19").Value = "  +2.71%  "

# Row 20
$ws.Range("D20").Value = "'14.34"
$ws.Range("E20").Value = "  +1.62%  "

# Row 21
$ws.Range("D21").Value = "'380.96"
$ws.Range("E21").Value = "  -0.46%  "

# Row 22
$ws.Range("D22").Value = "'8.16"
$ws.Range("E22").Value = "  +1.11%  "

# Row 23
$ws.Range("D23").Value = "'0.555"
$ws.Range("E23").Value = "  +3.76%  "

# Row 24
$ws.Range("D24").Value = "'72.97"
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("D26").Value = "'0.0000120"
$ws.Range("E26").Value = "  +1.51%  "

# Row 27
$ws.Range("D27").Value = "'10.04"
$ws.Range("E27").Value = "  +5.15%  "

# Row 28
$ws.Range("D28").Value = "'0.177"
$ws.Range("E28").Value = "  -0.47%  "

# Row 29
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("D30").Value = "'1.53"
$ws.Range("E30").Value = "  +10.80%  "

# Row 31
$ws.Range("D31").Value = "'6.17"
$ws.Range("E31").Value = "  +1.57%  "

# Row 32
$ws.Range("E32").Value = "  +2.76%  "

# Row 33
$ws.Range("D33").Value = "'23.60"
$ws.Range("E33").Value = "  +1.22%  "

# Row 34
$ws.Range("D34").Value = "'7.27"
$ws.Range("E34").Value = "  +5.98%  "

# Row 35
$ws.Range("E35").Value = "  +11.05%  "

# Row 36
$ws.Range("D36").Value = "'161.12"
$ws.Range("E36").Value = "  +1.15%  "

# Row 37
$ws.Range("D37").Value = "'1.94"
$ws.Range("E37").Value = "  +5.77%  "

# Row 38
$ws.Range("D38").Value = "'0.0778"
$ws.Range("E38").Value = "  +2.47%  "

# Row 39
$ws.Range("D39").Value = "2.912.78"
$ws.Range("E39").Value = "  +0.31%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'6.81"
$ws.Range("E40").Value = "  +6.25%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'4.66"
$ws.Range("E41").Value = "  +6.50%  "

# Row 42
$ws.Range("D42").Value = "'26.55"
$ws.Range("E42").Value = "  -0.57%  "

# Row 43
$ws.Range("D43").Value = "'0.0320"
$ws.Range("E43").Value = "  +1.17%  "

# Row 44
$ws.Range("D44").Value = "'43.01"
$ws.Range("E44").Value = "  +1.98%  "

# Row 45
$ws.Range("D45").Value = "'0.779"
$ws.Range("E45").Value = "  +3.40%  "

# Row 46
$ws.Range("D46").Value = "'25.99"
$ws.Range("E46").Value = "  +11.94%  "

# Row 47
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "'322.05"
$ws.Range("E47").Value = "  +11.24%  "

# Row 48
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.09"
$ws.Range("E48").Value = "  +2.95%  "

# Row 49
$ws.Range("E49").Value = "  +3.57%  "

# Row 50
$ws.Range("D50").Value = "'0.875"
$ws.Range("E50").Value = "  +4.31%  "

# Row 51
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'2.19"
$ws.Range("E51").Value = "  +0.17%  "
